$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently ends with:
#   ... row 318 (2024-05) / row 319 (2024-06) / row 320 ("Source..." footnote row)
# The update (UK EPU data refreshed to 2024q3) revises the May/June 2024 values
# and appends six new monthly rows for 2024 (Jul-Dec) before the footnote row,
# which is pushed from row 320 down to row 326.

# 1) Insert 6 blank rows right before the footnote row so it moves from 320 -> 326.
$ws.Range("A320:A325").EntireRow.Insert()

# 2) Revise the existing May/June 2024 index values.
$ws.Range("C318").Value = 144.78328886819094
$ws.Range("C319").Value = 146.59757407655673

# 3) Populate the 6 new rows for 2024 (Jul-Dec).
#    Column A must hold the same shared string already used for "2024" (used by
#    rows 314-319). Copy/PasteSpecial(values) from an existing "2024" cell keeps
#    it a shared-string text cell instead of turning it into a numeric literal.
$months = @(7, 8, 9, 10, 11, 12)
$values = @(183.28213802345289, 116.87180249328715, 184.40996223935207, 273.85700012719877, 327.23352118076377, 223.55022080880141)

$row = 320
for ($i = 0; $i -lt 6; $i++) {
  $ws.Range("A314").Copy()
  $ws.Cells.Item($row, 1).PasteSpecial(-4163)
  $ws.Cells.Item($row, 2).Value = $months[$i]
  $ws.Cells.Item($row, 3).Value = $values[$i]
  $row++
}

$excel.CutCopyMode = 0
